# Weekly update: a new Perejil price observation was recorded and inserted
# into the dataset at row 210, pushing every subsequent row down by one
# (the sheet appears to be ordered with the newest addition inserted in
# place rather than appended at the very end).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 210 - this shifts rows
# 210..312 down to 211..313 and extends the used range to row 313.
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new observation. Most
# columns repeat the constant "template" values shared by every data row
# in this sheet; only D (Fecha), J (Volumen), K/L/M (prices) and P
# (Precio $/Kg) are specific to this new record.
$ws.Cells.Item(210, 1).Value = 9
$ws.Cells.Item(210, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(210, 3).Value = "Metropolitana"
$ws.Cells.Item(210, 4).Value = 44609
$ws.Cells.Item(210, 5).Value = 13
$ws.Cells.Item(210, 6).Value = 100112044
$ws.Cells.Item(210, 7).Value = "Perejil"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 61
$ws.Cells.Item(210, 11).Value = 16000
$ws.Cells.Item(210, 12).Value = 18000
$ws.Cells.Item(210, 13).Value = 17016
$ws.Cells.Item(210, 14).Value = "$/docena de atados"
$ws.Cells.Item(210, 15).Value = "Región Metropolitana"
$ws.Cells.Item(210, 16).Value = 5672
$ws.Cells.Item(210, 17).Value = 3
$ws.Cells.Item(210, 18).Value = "Hortaliza"
